$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C4").AddComment("Thomas Smyth:`nThis one is a number so we are sure it can handle numbers.")
Write-Host "Comment added"
